# Regenerate the handback status report with new file GUIDs / xlf hash /
# timestamps (new handoff cycle). Mirrors the author's commit
# "Generate Report for Handback": two source files that used to be named
#   db09392f-62af-4bb1-863c-3b8e7e307202.md
#   ed8615df-9c41-4704-8c2b-f6ad065224d6.md
# are now named
#   f7179365-6677-4d53-b7bd-0ef081bd0637.md
#   ffffffda1da4-4572-4a73-8fee-606b23802e82.md
# with refreshed handoff/handback timestamps and a new xlf content hash.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "db09392f-62af-4bb1-863c-3b8e7e307202"
$oldGuid2 = "ed8615df-9c41-4704-8c2b-f6ad065224d6"
$newGuid1 = "f7179365-6677-4d53-b7bd-0ef081bd0637"
$newGuid2 = "ffffffda1da4-4572-4a73-8fee-606b23802e82"

$oldXlf1 = "7c5deeb5f46d211ec2ba4b56221f9575059a00d0"
$newXlf1 = "bb83099b247b3591770e4a3634e666fb1d9145c1"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A2").Value = "$newGuid1.md"
$wsOv.Range("G2").Value = "2016-08-28 19:00:25"

$wsOv.Range("A3").Value = "$newGuid2.md"
$wsOv.Range("G3").Value = "2016-08-28 19:00:25"

# hyperlinks on B2 / B3 -- the engine's Hyperlink property setters always
# append a fresh entry instead of updating in place, and Range.Hyperlinks.Delete()
# actually clears every hyperlink on the sheet -- so clear once and rebuild all
# of them with their original (unchanged) target URLs plus the new display text.
$ovB2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c39a9c38cd74d9ed3ce5f7faf14f3b3766e6af/e2e/$oldGuid1.md"
$ovB3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c39a9c38cd74d9ed3ce5f7faf14f3b3766e6af/e2e/$oldGuid2.md"

$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), $ovB2Target, "", "", "e2e\$newGuid1.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $ovB3Target, "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("G2").Value = "$newGuid1.$newXlf1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-28 19:00:20"
$wsZh.Range("J2").Value = "$newGuid1.$newXlf1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-28 19:00:36"

$wsZh.Range("G3").Value = "$newGuid1.$newXlf1.de-de.xlf"
$wsZh.Range("H3").Value = "2016-08-28 19:00:20"
$wsZh.Range("J3").Value = "$newGuid1.$newXlf1.de-de.xlf"
$wsZh.Range("K3").Value = "2016-08-28 19:00:36"

$zhA2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c39a9c38cd74d9ed3ce5f7faf14f3b3766e6af/e2e/$oldGuid1.md"
$zhI2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ecdd131ac5294e76dd14d5d07f6b1ea66aeb3d50/e2e/$oldGuid1.md"
$zhA3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c39a9c38cd74d9ed3ce5f7faf14f3b3766e6af/e2e/$oldGuid2.md"
$zhI3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ecdd131ac5294e76dd14d5d07f6b1ea66aeb3d50/e2e/$oldGuid2.md"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Target, "", "", "$newGuid1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhI2Target, "", "", "$newGuid1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Target, "", "", "$newGuid2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhI3Target, "", "", "$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("G2").Value = "$newGuid1.$newXlf1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-28 19:00:25"
$wsDe.Range("J2").Value = "$newGuid1.$newXlf1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-28 19:00:44"

$wsDe.Range("G3").Value = "$newGuid1.$newXlf1.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-28 19:00:25"
$wsDe.Range("J3").Value = "$newGuid1.$newXlf1.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-28 19:00:44"

$deA2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c39a9c38cd74d9ed3ce5f7faf14f3b3766e6af/e2e/$oldGuid1.md"
$deI2Target = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/467f9c6fdce5b23d00cdd6ba74497b8e4794d62c/e2e/$oldGuid1.md"
$deA3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c39a9c38cd74d9ed3ce5f7faf14f3b3766e6af/e2e/$oldGuid2.md"
$deI3Target = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/467f9c6fdce5b23d00cdd6ba74497b8e4794d62c/e2e/$oldGuid2.md"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Target, "", "", "$newGuid1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deI2Target, "", "", "$newGuid1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Target, "", "", "$newGuid2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deI3Target, "", "", "$newGuid2.md")

"Handback report regenerated"
